$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "2016 Rising Star Award in IBM cloud data server organization.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2016 Rising Star Award in IBM cloud data server organization.(3/500)",
    2)

$d.Content.Find.Execute(
    "2014 IBM IM Org Best Developer Award",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2014 IBM IM Org Best Developer Award.(1/200)",
    2)

$d.Content.Find.Execute(
    "2013 IBM CDL Technical Contest First Place Award(with Eric)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2013 IBM CDL Technical Contest First Place Award(co-auther)(1/3000)",
    2)
